$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

$lo = $ws.ListObjects.Item("Tabela1")
$newRow = $lo.ListRows.Add()

$ws.Range("A72").Value = 43972
$ws.Range("B72").Value = 73742
$ws.Range("C72").Value = 882
$ws.Range("D72").Value = 1468
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 21
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 106
$ws.Range("J72").Value = 0
